$d = $word.ActiveDocument

$replacements = @(
    @("10 ± 41", "10 ± 41 (2.6)"),
    @("182 ± 321", "182 ± 321 (98)"),
    @("35 ± 280", "35 ± 280 (7.8)"),
    @("170 ± 297", "170 ± 297 (93)"),
    @("9.3 ± 27", "9.3 ± 27 (3.2)"),
    @("197 ± 363", "197 ± 363 (112)"),
    @("22 ± 23", "22 ± 23 (17)"),
    @("106 ± 67", "106 ± 67 (92)"),
    @("1.7 ± 7.3", "1.7 ± 7.3 (0.49)"),
    @("174 ± 328", "174 ± 328 (96)"),
    @("1.4 ± 3.7", "1.4 ± 3.7 (0.5)"),
    @("123 ± 194", "123 ± 194 (75)"),
    @("18 ± 21", "18 ± 21 (13)"),
    @("100 ± 63", "100 ± 63 (88)"),
    @("1 ± 2.6", "1 ± 2.6 (0.4)"),
    @("311 ± 502", "311 ± 502 (181)"),
    @("2 ± 5.2", "2 ± 5.2 (0.72)"),
    @("74 ± 126", "74 ± 126 (42)"),
    @("8.6 ± 39", "8.6 ± 39 (2.2)"),
    @("141 ± 267", "141 ± 267 (74)"),
    @("1.8 ± 4.4", "1.8 ± 4.4 (0.68)"),
    @("99 ± 161", "99 ± 161 (59)"),
    @("3.2 ± 11", "3.2 ± 11 (0.85)"),
    @("149 ± 301", "149 ± 301 (80)"),
    @("8 ± 23", "8 ± 23 (2.8)"),
    @("89 ± 176", "89 ± 176 (39)"),
    @("28 ± 29", "28 ± 29 (22)"),
    @("93 ± 58", "93 ± 58 (82)"),
    @("1.9 ± 4.4", "1.9 ± 4.4 (0.71)"),
    @("142 ± 226", "142 ± 226 (88)"),
    @("27 ± 30", "27 ± 30 (20)"),
    @("92 ± 59", "92 ± 59 (81)"),
    @("9.7 ± 39", "9.7 ± 39 (2.9)"),
    @("195 ± 352", "195 ± 352 (101)"),
    @("9.8 ± 38", "9.8 ± 38 (3)"),
    @("98 ± 177", "98 ± 177 (51)")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
